$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Error codes")

# The "Error codes" table (A: numeric code, B: message) gains two new rows:
#   code 32 "Failed to open landscape costs file in readLandChange()"
#   code 38 "Invalid SMS cost read from costs file in readLandChange()"
# Insert blank rows at their final positions (bottom-most first so earlier
# row indices used below stay valid), which also pushes the existing
# formatting of the rows below (incl. the bold code column for the last
# ten rows) down by two rows, matching the target layout.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(12).Insert()

# Re-assert every row's final (code, message) pair in the now-correct order.
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Failed to open landscape habitat file in readLandscape()"
$ws.Cells.Item(3, 1).Value = 12
$ws.Cells.Item(3, 2).Value = "Failed to open landscape patch file in readLandscape()"
$ws.Cells.Item(4, 1).Value = 13
$ws.Cells.Item(4, 2).Value = "Invalid habitat number read from habitat file in readLandscape()"
$ws.Cells.Item(5, 1).Value = 14
$ws.Cells.Item(5, 2).Value = "Invalid patch number read from patch file in readLandscape()"
$ws.Cells.Item(6, 1).Value = 17
$ws.Cells.Item(6, 2).Value = "Invalid quality score read from habitat file in readLandscape()"
$ws.Cells.Item(7, 1).Value = 19
$ws.Cells.Item(7, 2).Value = "Invalid no. of habitats specified"
$ws.Cells.Item(8, 1).Value = 21
$ws.Cells.Item(8, 2).Value = "Failed to open initial distribution file in readDistribution()"
$ws.Cells.Item(9, 1).Value = 22
$ws.Cells.Item(9, 2).Value = "Invalid value read from initial distribution file in readDistribution()"
$ws.Cells.Item(10, 1).Value = 30
$ws.Cells.Item(10, 2).Value = "Failed to open landscape habitat file in readLandChange()"
$ws.Cells.Item(11, 1).Value = 31
$ws.Cells.Item(11, 2).Value = "Failed to open landscape patch file in readLandChange()"
$ws.Cells.Item(12, 1).Value = 32
$ws.Cells.Item(12, 2).Value = "Failed to open landscape costs file in readLandChange()"
$ws.Cells.Item(13, 1).Value = 33
$ws.Cells.Item(13, 2).Value = "Invalid habitat number read from habitat file in readLandChange()"
$ws.Cells.Item(14, 1).Value = 34
$ws.Cells.Item(14, 2).Value = "Invalid patch number read from habitat file in readLandChange()"
$ws.Cells.Item(15, 1).Value = 36
$ws.Cells.Item(15, 2).Value = "Invalid 'no data' cell read from habitat file in readLandChange()"
$ws.Cells.Item(16, 1).Value = 37
$ws.Cells.Item(16, 2).Value = "Invalid quality score read from habitat file in readLandChange()"
$ws.Cells.Item(17, 1).Value = 38
$ws.Cells.Item(17, 2).Value = "Invalid SMS cost read from costs file in readLandChange()"
$ws.Cells.Item(18, 1).Value = 41
$ws.Cells.Item(18, 2).Value = "Failed to open dynamic landscape file in ReadDynLandFile()"
$ws.Cells.Item(19, 1).Value = 51
$ws.Cells.Item(19, 2).Value = "Invalid header in costs file"
$ws.Cells.Item(20, 1).Value = 52
$ws.Cells.Item(20, 2).Value = "Dimensions of cost map file differ from dimensions of landscape"
$ws.Cells.Item(21, 1).Value = 53
$ws.Cells.Item(21, 2).Value = "Origin of cost map file differs from origin of landscape"
$ws.Cells.Item(22, 1).Value = 54
$ws.Cells.Item(22, 2).Value = "Invalid value in costs file"
$ws.Cells.Item(23, 1).Value = 101
$ws.Cells.Item(23, 2).Value = "Error in the parameters file: local environmental stochasticity is not allowed with a patch-based model"
$ws.Cells.Item(24, 1).Value = 102
$ws.Cells.Item(24, 2).Value = "Error in the parameters file: local extinction probability is not allowed with a patch-based model"
$ws.Cells.Item(25, 1).Value = 103
$ws.Cells.Item(25, 2).Value = "Error in the parameters file: the occupancy output is possible only with more than 1 replicate"
$ws.Cells.Item(26, 1).Value = 104
$ws.Cells.Item(26, 2).Value = "Error in the parameters file: the output traits by row is not allowed for patch-based models"
$ws.Cells.Item(27, 1).Value = 105
$ws.Cells.Item(27, 2).Value = "Error in the parameters file: the connectivity matrix can be computed only for patch-based models"
$ws.Cells.Item(28, 1).Value = 300
$ws.Cells.Item(28, 2).Value = "Simulation mis-match within ReadEmigration()"
$ws.Cells.Item(29, 1).Value = 301
$ws.Cells.Item(29, 2).Value = "Error in the emigration file: sex-dependent emigration is not possible with asexual models"
$ws.Cells.Item(30, 1).Value = 303
$ws.Cells.Item(30, 2).Value = "Error in the emigration file: stage-dependent emigration is not possible without stage structure"
$ws.Cells.Item(31, 1).Value = 400
$ws.Cells.Item(31, 2).Value = "Simulation mis-match within ReadTransfer()"
$ws.Cells.Item(32, 1).Value = 401
$ws.Cells.Item(32, 2).Value = "Error in the transfer file: sex-dependent kernels are not possible with asexual models"
$ws.Cells.Item(33, 1).Value = 403
$ws.Cells.Item(33, 2).Value = "Error in the transfer file: stage-dependent kernels are not possible without stage structure"
$ws.Cells.Item(34, 1).Value = 434
$ws.Cells.Item(34, 2).Value = "Error in the transfer files: with the set type of landscape the per-step mortality cannot be habitat-dependent"
$ws.Cells.Item(35, 1).Value = 440
$ws.Cells.Item(35, 2).Value = "Internal error in ReadTransfer(): unknown TransferType"
$ws.Cells.Item(36, 1).Value = 500
$ws.Cells.Item(36, 2).Value = "Simulation mis-match within ReadSettlement()"
$ws.Cells.Item(37, 1).Value = 501
$ws.Cells.Item(37, 2).Value = "Error in the settlement file: sex-dependent settlement is not possible with asexual models"
$ws.Cells.Item(38, 1).Value = 502
$ws.Cells.Item(38, 2).Value = "Error in the settlement file: stage-dependent settlement is not possible without stage structure"
$ws.Cells.Item(39, 1).Value = 503
$ws.Cells.Item(39, 2).Value = "Error in the settlement file: the chosen settlement rule is possible only with overlapping generations"
$ws.Cells.Item(40, 1).Value = 504
$ws.Cells.Item(40, 2).Value = "Error in the settlement file: the rule 'find a mate' is possible only with sexual models"
$ws.Cells.Item(41, 1).Value = 505
$ws.Cells.Item(41, 2).Value = "Error in the settlement file: the settlement rule chosen is possible only with overlapping generations"
$ws.Cells.Item(42, 1).Value = 507
$ws.Cells.Item(42, 2).Value = "Error in the settlement file: the rule 'find a mate' is possible only with sexual models"
$ws.Cells.Item(43, 1).Value = 508
$ws.Cells.Item(43, 2).Value = "Error in the settlement file: sex-dependent settlement is not possible with asexual models"
$ws.Cells.Item(44, 1).Value = 509
$ws.Cells.Item(44, 2).Value = "Error in the settlement file: stage-dependent settlement is not possible without stage structure"
$ws.Cells.Item(45, 1).Value = 510
$ws.Cells.Item(45, 2).Value = "Error in the settlement file: the rule 'find a mate' is possible only with sexual models"
$ws.Cells.Item(46, 1).Value = 601
$ws.Cells.Item(46, 2).Value = "Error in the initialisation parameters file: the initialisation cannot be from species distribution when no species distribution has been loaded"
$ws.Cells.Item(47, 1).Value = 602
$ws.Cells.Item(47, 2).Value = "Error in the initialisation parameters file: maximum number of initial cells exceeded"
$ws.Cells.Item(48, 1).Value = 603
$ws.Cells.Item(48, 2).Value = "Error in the initialisation parameters file: initial limits exceed landscape dimensions"
$ws.Cells.Item(49, 1).Value = 901
$ws.Cells.Item(49, 2).Value = "Error in the landscape file: Invalid dimensions. ensure Y >= X"
$ws.Cells.Item(50, 1).Value = 902
$ws.Cells.Item(50, 2).Value = "Error in the landscape file: Invalid dimensions. ensure X & Y are powers of 2 + 1."
